# Apply the update that happened on 2024-09-05 10:57:15 to the "2024" sheet.
# A new September transaction ("bal axisbank axis") was recorded, which
# pushes every existing row of the "Others" category (and everything that
# follows it, including the "hdfc"/August rows and the "Broadband" label)
# down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new blank row above row 31; this shifts rows 31:73 down to 32:74
# (dimension grows from A1:Y73 to A1:Y74) exactly like Excel does when a
# new entry is inserted above the most-recent September transaction.
$ws.Rows.Item(31).Insert()

# Populate the freshly inserted row with the new September transaction.
$ws.Range("R31").Value = "bal axisbank axis"
$ws.Range("S31").Value = "2024-09-05 16:26:56"
